$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.974.33"
$ws.Range("E2").Value = "  +5.33%  "

$ws.Range("D3").Value = "2.417.70"
$ws.Range("E3").Value = "  +2.00%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").Value = "'573.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.43%  "

$ws.Range("D6").Value = "'145.48"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.42%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("E8").Value = "  +2.62%  "

$ws.Range("D9").Value = "2.449.11"
$ws.Range("E9").Value = "  +3.55%  "

$ws.Range("E10").Value = "  +5.87%  "

$ws.Range("D11").Value = "'0.160"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").Value = "'5.23"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.93%  "

$ws.Range("E13").Value = "  +4.78%  "

$ws.Range("D14").Value = "'27.29"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.85%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000178"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.90%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.868.56"
$ws.Range("E16").Value = "  +2.22%  "

$ws.Range("D17").Value = "62.715.00"
$ws.Range("E17").Value = "  +4.93%  "

$ws.Range("D18").Value = "2.449.03"
$ws.Range("E18").Value = "  +2.89%  "

$ws.Range("D19").Value = "'7.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.27%  "

$ws.Range("D20").Value = "'10.97"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.01%  "

$ws.Range("D21").Value = "'328.63"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.85%  "

$ws.Range("E22").Value = "  +2.27%  "

$ws.Range("D23").Value = "'2.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +12.13%  "

$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("E25").Value = "  +2.24%  "

$ws.Range("D26").Value = "'628.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +11.57%  "

$ws.Range("E27").Value = "  +10.62%  "

$ws.Range("D28").Value = "'8.46"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.28%  "

$ws.Range("D29").Value = "0.0₃0982"
$ws.Range("E29").Value = "  +6.21%  "

$ws.Range("D30").Value = "2.532.99"
$ws.Range("E30").Value = "  +2.00%  "

$ws.Range("D31").Value = "'8.19"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.53%  "

$ws.Range("E32").Value = "  +8.70%  "

$ws.Range("D33").Value = "'0.138"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.64%  "

$ws.Range("E34").Value = "  +3.85%  "

$ws.Range("D35").Value = "'1.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.72%  "

$ws.Range("D36").Value = "'0.995"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  +4.85%  "

$ws.Range("E38").Value = "  +2.08%  "

$ws.Range("D39").Value = "'151.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.76%  "

$ws.Range("D40").Value = "'5.39"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.31%  "

$ws.Range("D41").Value = "'18.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.99%  "

$ws.Range("E42").Value = "  +14.23%  "

$ws.Range("E43").Value = "  +7.47%  "

$ws.Range("D45").Value = "0.0₆0302"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").Value = "'145.01"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.41%  "

$ws.Range("E47").Value = "  +2.08%  "

$ws.Range("D48").Value = "'20.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.23%  "

$ws.Range("D49").Value = "'0.602"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.18%  "

$ws.Range("D50").Value = "'0.0516"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.35%  "

$ws.Range("E51").Value = "  +2.71%  "

